$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15.11121702360658
$ws.Range("C2").Value = 11.36169486345986
$ws.Range("D2").Value = 5.957643046581369
$ws.Range("E2").Value = 16.59538305590427
$ws.Range("G2").Value = 19.89179704584907
$ws.Range("H2").Value = 11.78354802379599
$ws.Range("O2").Value = 16.77201324948246

$ws.Range("B3").Value = 14.27686674066822
$ws.Range("C3").Value = 10.75757236940133
$ws.Range("D3").Value = 5.833783793116681
$ws.Range("E3").Value = 15.64683130777179
$ws.Range("G3").Value = 19.84163559637258
$ws.Range("H3").Value = 11.84814371662591
$ws.Range("O3").Value = 16.85584347804583

$ws.Range("B4").Value = 13.73858331108573
$ws.Range("C4").Value = 10.36741900117662
$ws.Range("D4").Value = 5.758206947671873
$ws.Range("E4").Value = 15.03910494268981
$ws.Range("G4").Value = 19.8258643230265
$ws.Range("H4").Value = 11.89123420485735
$ws.Range("O4").Value = 16.91464351562693

$ws.Range("B5").Value = 13.51286102365914
$ws.Range("C5").Value = 10.20370216282162
$ws.Range("D5").Value = 5.727576091240363
$ws.Range("E5").Value = 14.78535355371565
$ws.Range("G5").Value = 19.82319804593656
$ws.Range("H5").Value = 11.90965226478539
$ws.Range("O5").Value = 16.94043110117973

$ws.Range("B6").Value = 13.47500102895907
$ws.Range("C6").Value = 10.17623521286279
$ws.Range("D6").Value = 5.722501390816414
$ws.Range("E6").Value = 14.74285887715245
$ws.Range("G6").Value = 19.82298173485711
$ws.Range("H6").Value = 11.91276230090092
$ws.Range("O6").Value = 16.94482290852612

$ws.Range("B7").Value = 13.73556467559954
$ws.Range("C7").Value = 10.36523003754988
$ws.Range("D7").Value = 5.757793107425508
$ws.Range("E7").Value = 15.03570706133702
$ws.Range("G7").Value = 19.82581316898048
$ws.Range("H7").Value = 11.89147912691567
$ws.Range("O7").Value = 16.91498392456689

$ws.Range("B8").Value = 14.82901732974186
$ws.Range("C8").Value = 11.15744003144679
$ws.Range("D8").Value = 5.914866919145019
$ws.Range("E8").Value = 16.27372027605158
$ws.Range("G8").Value = 19.87137413819626
$ws.Range("H8").Value = 11.80510676086807
$ws.Range("O8").Value = 16.79938716703175

$ws.Range("B9").Value = 16.76128862464084
$ws.Range("C9").Value = 12.55477626120955
$ws.Range("D9").Value = 6.224565504079205
$ws.Range("E9").Value = 18.6098675262162
$ws.Range("G9").Value = 20.08030674715361
$ws.Range("H9").Value = 11.66310405576438
$ws.Range("O9").Value = 16.63156430761904

$ws.Range("B10").Value = 18.0457774974455
$ws.Range("C10").Value = 13.48240434197996
$ws.Range("D10").Value = 6.450437955877692
$ws.Range("E10").Value = 20.26729824874777
$ws.Range("G10").Value = 20.30651053177734
$ws.Range("H10").Value = 11.57569581534308
$ws.Range("O10").Value = 16.54509321981104

$ws.Range("B11").Value = 18.59994177894406
$ws.Range("C11").Value = 13.88237657240083
$ws.Range("D11").Value = 6.552300009501745
$ws.Range("E11").Value = 20.97889756704681
$ws.Range("G11").Value = 20.42496068047622
$ws.Range("H11").Value = 11.53965645673632
$ws.Range("O11").Value = 16.51394121311515

$ws.Range("B12").Value = 18.80540317361982
$ws.Range("C12").Value = 14.03063877644278
$ws.Range("D12").Value = 6.590703467256724
$ws.Range("E12").Value = 21.24230585423862
$ws.Range("G12").Value = 20.47201694699935
$ws.Range("H12").Value = 11.52654891497538
$ws.Range("O12").Value = 16.5033360403837

$ws.Range("B13").Value = 18.76134926852175
$ws.Range("C13").Value = 13.99885056470708
$ws.Range("D13").Value = 6.582440771614705
$ws.Range("E13").Value = 21.18584509324717
$ws.Range("G13").Value = 20.46178524771847
$ws.Range("H13").Value = 11.52934778139948
$ws.Range("O13").Value = 16.50556684279704

$ws.Range("B14").Value = 18.6169334198846
$ws.Range("C14").Value = 13.89463846106364
$ws.Range("D14").Value = 6.555463112118518
$ws.Range("E14").Value = 21.0006895589115
$ws.Range("G14").Value = 20.42878812695709
$ws.Range("H14").Value = 11.53856724808841
$ws.Range("O14").Value = 16.51304476106891

$ws.Range("B15").Value = 18.52790168483289
$ws.Range("C15").Value = 13.83038813082581
$ws.Range("D15").Value = 6.538915247849139
$ws.Range("E15").Value = 20.88648841439135
$ws.Range("G15").Value = 20.40886201378709
$ws.Range("H15").Value = 11.54428486406438
$ws.Range("O15").Value = 16.51778077699928

$ws.Range("B16").Value = 18.00894998564705
$ws.Range("C16").Value = 13.45581925354399
$ws.Range("D16").Value = 6.443759529323128
$ws.Range("E16").Value = 20.21994383027322
$ws.Range("G16").Value = 20.2990796969794
$ws.Range("H16").Value = 11.57812636651473
$ws.Range("O16").Value = 16.54729502475939

$ws.Range("B17").Value = 17.68282492605185
$ws.Range("C17").Value = 13.2203694278914
$ws.Range("D17").Value = 6.385126551497812
$ws.Range("E17").Value = 19.80021182993872
$ws.Range("G17").Value = 20.23569270271395
$ws.Range("H17").Value = 11.59984373419694
$ws.Range("O17").Value = 16.56750782151632

$ws.Range("B18").Value = 17.49241018926383
$ws.Range("C18").Value = 13.08287460684205
$ws.Range("D18").Value = 6.351321131785789
$ws.Range("E18").Value = 19.55480257614744
$ws.Range("G18").Value = 20.2007003282147
$ws.Range("H18").Value = 11.61268515997848
$ws.Range("O18").Value = 16.57990319614915

$ws.Range("B19").Value = 17.42745370882009
$ws.Range("C19").Value = 13.03596680483899
$ws.Range("D19").Value = 6.339862549235982
$ws.Range("E19").Value = 19.47102479350259
$ws.Range("G19").Value = 20.1891052227211
$ws.Range("H19").Value = 11.61709305413919
$ws.Range("O19").Value = 16.58423176400697

$ws.Range("B20").Value = 17.71783545253241
$ws.Range("C20").Value = 13.24564800178898
$ws.Range("D20").Value = 6.391376854571843
$ws.Range("E20").Value = 19.84530576951444
$ws.Range("G20").Value = 20.24228883613321
$ws.Range("H20").Value = 11.59749561029538
$ws.Range("O20").Value = 16.56527638578966

$ws.Range("B21").Value = 18.65947131426379
$ws.Range("C21").Value = 13.92533511808278
$ws.Range("D21").Value = 6.563392024834252
$ws.Range("E21").Value = 21.0552383919055
$ws.Range("G21").Value = 20.43842074012616
$ws.Range("H21").Value = 11.53584458255239
$ws.Range("O21").Value = 16.51081586871878

$ws.Range("B22").Value = 19.24928888762082
$ws.Range("C22").Value = 14.35089594716779
$ws.Range("D22").Value = 6.674808380711819
$ws.Range("E22").Value = 21.81069904208332
$ws.Range("G22").Value = 20.57941615065623
$ws.Range("H22").Value = 11.49870038241534
$ws.Range("O22").Value = 16.48217340032562

$ws.Range("B23").Value = 18.93684827652883
$ws.Range("C23").Value = 14.12548200584514
$ws.Range("D23").Value = 6.615448516415749
$ws.Range("E23").Value = 21.41071413150972
$ws.Range("G23").Value = 20.5030054055013
$ws.Range("H23").Value = 11.5182354117503
$ws.Range("O23").Value = 16.49681985532698

$ws.Range("B24").Value = 17.70201628800093
$ws.Range("C24").Value = 13.23422619535557
$ws.Range("D24").Value = 6.388551390068161
$ws.Range("E24").Value = 19.82493157216451
$ws.Range("G24").Value = 20.23930220623176
$ws.Range("H24").Value = 11.59855608907689
$ws.Range("O24").Value = 16.56628280565215

$ws.Range("B25").Value = 16.26194994408302
$ws.Range("C25").Value = 12.19392569039743
$ws.Range("D25").Value = 6.140888214780666
$ws.Range("E25").Value = 17.96169822624853
$ws.Range("G25").Value = 20.01096595254244
$ws.Range("H25").Value = 11.69856502193227
$ws.Range("O25").Value = 16.67056337191977
